$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2513.9092
$ws.Range("I100").Value = 2842.1667
$ws.Range("J100").Value = 2120
$ws.Range("K100").Value = 2842.1667
$ws.Range("L100").Value = 2120
$ws.Range("M100").Value = -2301.1667
$ws.Range("N100").Value = -3202
$ws.Range("H106").Value = 2231.3333
$ws.Range("J106").Value = 2353
$ws.Range("L106").Value = 2353
$ws.Range("N106").Value = -3615
$ws.Range("H137").Value = 1704.7667
$ws.Range("I137").Value = 964.7222
$ws.Range("K137").Value = 2894.1666
$ws.Range("M137").Value = -344.1666
$ws.Range("H138").Value = 4720.48
$ws.Range("I138").Value = 1966.3334
$ws.Range("J138").Value = 5590.2104
$ws.Range("K138").Value = 5899.0002
$ws.Range("L138").Value = 16770.6312
$ws.Range("M138").Value = -759.0002000000004
$ws.Range("N138").Value = -27050.6312

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3083.9333
$ws.Range("I61").Value = 3018.5
$ws.Range("K61").Value = 3018.5
$ws.Range("M61").Value = -2806.5
$ws.Range("H74").Value = 2283.6667
$ws.Range("J74").Value = 4204.5
$ws.Range("L74").Value = 4204.5
$ws.Range("N74").Value = -5952.5
$ws.Range("H77").Value = 2283.6667
$ws.Range("J77").Value = 4204.5
$ws.Range("L77").Value = 21022.5
$ws.Range("N77").Value = -29758.5
$ws.Range("H102").Value = 4263.0527
$ws.Range("I102").Value = 2849.8333
$ws.Range("J102").Value = 6685.7144
$ws.Range("K102").Value = 2849.8333
$ws.Range("L102").Value = 6685.7144
$ws.Range("M102").Value = -1227.8333
$ws.Range("N102").Value = -9929.714400000001
$ws.Range("H132").Value = 1862.125
$ws.Range("I132").Value = 1708.3334
$ws.Range("J132").Value = 2323.5
$ws.Range("K132").Value = 5125.0002
$ws.Range("L132").Value = 6970.5
$ws.Range("M132").Value = -2595.0002
$ws.Range("N132").Value = -12030.5
$ws.Range("H136").Value = 3083.9333
$ws.Range("I136").Value = 3018.5
$ws.Range("K136").Value = 9055.5
$ws.Range("M136").Value = -6505.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2744.6
$ws.Range("I20").Value = 3378.7144
$ws.Range("K20").Value = 3378.7144
$ws.Range("M20").Value = -3131.7144
$ws.Range("H105").Value = 1002081.25
$ws.Range("I105").Value = 1540586.5
$ws.Range("K105").Value = 1540586.5
$ws.Range("M105").Value = -1538839.5
$ws.Range("H107").Value = 5711.4116
$ws.Range("I107").Value = 1701.875
$ws.Range("K107").Value = 1701.875
$ws.Range("M107").Value = 218.125
$ws.Range("H134").Value = 3721.0454
$ws.Range("I134").Value = 882.3158
$ws.Range("K134").Value = 2646.9474
$ws.Range("M134").Value = -111.9474

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5576.852
$ws.Range("I31").Value = 2958
$ws.Range("K31").Value = 2958
$ws.Range("M31").Value = -2663
$ws.Range("H34").Value = 5576.852
$ws.Range("I34").Value = 2958
$ws.Range("K34").Value = 2958
$ws.Range("M34").Value = -2756
$ws.Range("H105").Value = 3401.8
$ws.Range("I105").Value = 1499
$ws.Range("K105").Value = 1499
$ws.Range("M105").Value = 248
$ws.Range("H134").Value = 3980.4167
$ws.Range("I134").Value = 3897.9092
$ws.Range("K134").Value = 11693.7276
$ws.Range("M134").Value = -9158.7276

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 389.66666
$ws.Range("I97").Value = 420
$ws.Range("J97").Value = 374.5
$ws.Range("K97").Value = 1260
$ws.Range("L97").Value = 1123.5
$ws.Range("M97").Value = -764
$ws.Range("N97").Value = -2115.5
$ws.Range("H121").Value = 561.5714
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -5620
$ws.Range("H131").Value = 2238.3333
$ws.Range("I131").Value = 1832.3334
$ws.Range("J131").Value = 2373.6667
$ws.Range("K131").Value = 5497.0002
$ws.Range("L131").Value = 7121.000100000001
$ws.Range("M131").Value = -457.0002000000004
$ws.Range("N131").Value = -17201.0001
$ws.Range("H132").Value = 4812.25
$ws.Range("J132").Value = 4833.3335
$ws.Range("L132").Value = 43500.0015
$ws.Range("N132").Value = -48560.0015

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4750
$ws.Range("J29").Value = 4750
$ws.Range("L29").Value = 4750
$ws.Range("N29").Value = -5330
$ws.Range("H102").Value = 1688.75
$ws.Range("I102").Value = 941.6667
$ws.Range("J102").Value = 3930
$ws.Range("K102").Value = 941.6667
$ws.Range("L102").Value = 3930
$ws.Range("M102").Value = 680.3333
$ws.Range("N102").Value = -7174
$ws.Range("H132").Value = 102421.1
$ws.Range("I132").Value = 113466.336
$ws.Range("K132").Value = 340399.008
$ws.Range("M132").Value = -337869.008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6012.5625
$ws.Range("I46").Value = 2466.6667
$ws.Range("K46").Value = 2466.6667
$ws.Range("M46").Value = -2278.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 12300000
$ws.Range("J5").Value = 10500000
$ws.Range("L5").Value = 10500000
$ws.Range("N5").Value = -10500224
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51144
$ws.Range("H81").Value = 824.5
$ws.Range("I81").Value = 824.5
$ws.Range("K81").Value = 1649
$ws.Range("M81").Value = -588
$ws.Range("H84").Value = 824.5
$ws.Range("I84").Value = 824.5
$ws.Range("K84").Value = 8245
$ws.Range("M84").Value = -2941
$ws.Range("H100").Value = 2401.3333
$ws.Range("I100").Value = 2401.3333
$ws.Range("K100").Value = 4802.6666
$ws.Range("M100").Value = -4261.6666
$ws.Range("H107").Value = 923.5
$ws.Range("I107").Value = 981.3333
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 2943.9999
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = -1023.9999
$ws.Range("N107").Value = -6090
$ws.Range("H122").Value = 2759.2727
$ws.Range("I122").Value = 2589.4211
$ws.Range("K122").Value = 7768.263300000001
$ws.Range("M122").Value = -5318.263300000001
$ws.Range("H132").Value = 1757.4
$ws.Range("I132").Value = 1446.75
$ws.Range("K132").Value = 4340.25
$ws.Range("M132").Value = -1810.25
